$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '75.828.93'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '2.878.77'
$ws.Range('E3').Value = '  +5.21%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '195.16'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '597.08'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +1.44%  '
$ws.Range('E9').Value = '  -3.98%  '
$ws.Range('D10').Value = '2.878.83'
$ws.Range('E10').Value = '  +5.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.396'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +8.56%  '
$ws.Range('E12').Value = '  -1.61%  '
$ws.Range('E13').Value = '  +1.90%  '
$ws.Range('D14').Value = '3.413.45'
$ws.Range('E14').Value = '  +7.23%  '
$ws.Range('D15').Value = '75.752.22'
$ws.Range('E15').Value = '  +0.26%  '
$ws.Range('E16').Value = '  -1.79%  '
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').Value = '2.888.49'
$ws.Range('E18').Value = '  +7.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.89'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -7.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.53'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '379.80'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('E22').Value = '  -1.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.13'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('E24').Value = '  +0.59%  '
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('D26').Value = '3.042.10'
$ws.Range('E26').Value = '  +7.49%  '
$ws.Range('E27').Value = '  -1.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.73'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('E29').Value = '  +7.65%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('E31').Value = '  -2.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '506.75'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.63%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.77'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.81'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '164.49'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.18'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.113'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -6.51%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '182.90'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +5.39%  '
$ws.Range('E41').Value = '  -0.13%  '
$ws.Range('E42').Value = '  +2.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.99'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.81%  '
$ws.Range('E44').Value = '  -3.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0909'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +6.36%  '
$ws.Range('E46').Value = '  -1.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.41'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.09%  '
$ws.Range('E48').Value = '  -4.53%  '
$ws.Range('E49').Value = '  +4.38%  '
$ws.Range('E51').Value = '  +0.41%  '
